$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlink on C2 (mailto:mauruco009@gmail.com)
foreach ($h in $ws.Hyperlinks) {
    $h.Delete()
}

# Update row 2 data with the corrected client info
$ws.Range("A2").Value = "CONJUNTO DE PRUEBA 1"
$ws.Range("B2").Value = 111111
$ws.Range("C2").Value = "santigarcia2321@gmail.com"
$ws.Range("D2").Value = 111111

# Drop the trailing UID/ERROR columns (U:Z) that are no longer needed
$ws.Range("U1:Z2").ClearContents()

# Reflect the current selection left by the edit
$ws.Range("T1:Y1048576").Select()
